$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 54411.527
$ws.Range("I28").Value = 57385.055
$ws.Range("J28").Value = 888
$ws.Range("K28").Value = 57385.055
$ws.Range("L28").Value = 888
$ws.Range("M28").Value = -56900.055
$ws.Range("N28").Value = -1858
$ws.Range("H33").Value = 2050622.1
$ws.Range("I33").Value = 2158536.5
$ws.Range("K33").Value = 2158536.5
$ws.Range("M33").Value = -2158307.5
$ws.Range("H69").Value = 11610.125
$ws.Range("J69").Value = 11610.125
$ws.Range("L69").Value = 34830.375
$ws.Range("N69").Value = -36578.375
$ws.Range("H72").Value = 11610.125
$ws.Range("J72").Value = 11610.125
$ws.Range("L72").Value = 104491.125
$ws.Range("N72").Value = -113227.125
$ws.Range("H92").Value = 342.6842
$ws.Range("I92").Value = 310.0625
$ws.Range("K92").Value = 310.0625
$ws.Range("M92").Value = 937.9375
$ws.Range("H138").Value = 5473.229
$ws.Range("J138").Value = 6917.982
$ws.Range("L138").Value = 20753.946
$ws.Range("N138").Value = -31033.946

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 155.14285
$ws.Range("J4").Value = 95.5
$ws.Range("L4").Value = 95.5
$ws.Range("N4").Value = -327.5
$ws.Range("H28").Value = 7267.75
$ws.Range("I28").Value = 7267.75
$ws.Range("K28").Value = 7267.75
$ws.Range("M28").Value = -7075.75
$ws.Range("H61").Value = 4401.9756
$ws.Range("I61").Value = 4418.3438
$ws.Range("K61").Value = 4418.3438
$ws.Range("M61").Value = -4206.3438
$ws.Range("H74").Value = 3212.889
$ws.Range("I74").Value = 3679.3845
$ws.Range("K74").Value = 3679.3845
$ws.Range("M74").Value = -2805.3845
$ws.Range("H77").Value = 3212.889
$ws.Range("I77").Value = 3679.3845
$ws.Range("K77").Value = 18396.9225
$ws.Range("M77").Value = -14028.9225
$ws.Range("H97").Value = 648.58826
$ws.Range("I97").Value = 683.8214
$ws.Range("J97").Value = 484.16666
$ws.Range("K97").Value = 683.8214
$ws.Range("L97").Value = 484.16666
$ws.Range("M97").Value = -187.8214
$ws.Range("N97").Value = -1476.16666
$ws.Range("H99").Value = 7267.75
$ws.Range("I99").Value = 7267.75
$ws.Range("K99").Value = 7267.75
$ws.Range("M99").Value = -4272.75
$ws.Range("H122").Value = 4414.4814
$ws.Range("I122").Value = 2906.5715
$ws.Range("K122").Value = 8719.7145
$ws.Range("M122").Value = -6269.7145
$ws.Range("H132").Value = 4325.522
$ws.Range("I132").Value = 3269.932
$ws.Range("K132").Value = 9809.795999999998
$ws.Range("M132").Value = -7279.795999999998
$ws.Range("H136").Value = 4401.9756
$ws.Range("I136").Value = 4418.3438
$ws.Range("K136").Value = 13255.0314
$ws.Range("M136").Value = -10705.0314

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 20410706
$ws.Range("I20").Value = 26318180
$ws.Range("J20").Value = 3068.7273
$ws.Range("K20").Value = 26318180
$ws.Range("L20").Value = 3068.7273
$ws.Range("M20").Value = -26317933
$ws.Range("N20").Value = -3562.7273
$ws.Range("H86").Value = 812134.4
$ws.Range("I86").Value = 1309525
$ws.Range("K86").Value = 1309525
$ws.Range("M86").Value = -1308402
$ws.Range("H89").Value = 812134.4
$ws.Range("I89").Value = 1309525
$ws.Range("K89").Value = 6547625
$ws.Range("M89").Value = -6542009
$ws.Range("H94").Value = 1856.5161
$ws.Range("I94").Value = 1961.5385
$ws.Range("J94").Value = 1310.4
$ws.Range("K94").Value = 1961.5385
$ws.Range("L94").Value = 1310.4
$ws.Range("M94").Value = -1510.5385
$ws.Range("N94").Value = -2212.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 233.28
$ws.Range("I7").Value = 233.28
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 233.28
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -120.28
$ws.Range("H31").Value = 35288.133
$ws.Range("I31").Value = 1046.909
$ws.Range("J31").Value = 55112
$ws.Range("K31").Value = 1046.909
$ws.Range("L31").Value = 55112
$ws.Range("M31").Value = -751.9090000000001
$ws.Range("N31").Value = -55702
$ws.Range("H34").Value = 35288.133
$ws.Range("I34").Value = 1046.909
$ws.Range("J34").Value = 55112
$ws.Range("K34").Value = 1046.909
$ws.Range("L34").Value = 55112
$ws.Range("M34").Value = -844.9090000000001
$ws.Range("N34").Value = -55516
$ws.Range("H132").Value = 2432.1052
$ws.Range("I132").Value = 2078.8667
$ws.Range("K132").Value = 6236.6001
$ws.Range("M132").Value = -3706.6001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 11875466
$ws.Range("I11").Value = 18095492
$ws.Range("J11").Value = 868.9091
$ws.Range("K11").Value = 54286476
$ws.Range("L11").Value = 2606.7273
$ws.Range("M11").Value = -54286336
$ws.Range("N11").Value = -2886.7273
$ws.Range("H56").Value = 6480.909
$ws.Range("I56").Value = 6480.909
$ws.Range("K56").Value = 6480.909
$ws.Range("M56").Value = -5950.909
$ws.Range("H107").Value = 69367.3
$ws.Range("J107").Value = 94187.55
$ws.Range("L107").Value = 282562.65
$ws.Range("N107").Value = -286402.65
$ws.Range("H132").Value = 396876.2
$ws.Range("I132").Value = 93310.73
$ws.Range("K132").Value = 839796.5699999999
$ws.Range("M132").Value = -837266.5699999999
$ws.Range("H134").Value = 2138
$ws.Range("I134").Value = 2138
$ws.Range("K134").Value = 6414
$ws.Range("M134").Value = -1344

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 98500
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 98500
$ws.Range("K19").Value = 0
$ws.Range("L19").ClearContents()
$ws.Range("M19").Value = 98500
$ws.Range("N19").Value = -99076
$ws.Range("H70").Value = 20000.455
$ws.Range("I70").Value = 12856.143
$ws.Range("K70").Value = 12856.143
$ws.Range("M70").Value = -12586.143
$ws.Range("H73").Value = 20000.455
$ws.Range("I73").Value = 12856.143
$ws.Range("K73").Value = 12856.143
$ws.Range("M73").Value = -11920.143
$ws.Range("H102").Value = 2018.4286
$ws.Range("I102").Value = 1021.8276
$ws.Range("J102").Value = 4241.615
$ws.Range("K102").Value = 1021.8276
$ws.Range("L102").Value = 4241.615
$ws.Range("M102").Value = 600.1724
$ws.Range("N102").Value = -7485.615
$ws.Range("H126").Value = 4173.8184
$ws.Range("I126").Value = 3486.6667
$ws.Range("J126").Value = 4431.5
$ws.Range("K126").Value = 10460.0001
$ws.Range("L126").Value = 13294.5
$ws.Range("M126").Value = -7990.000100000001
$ws.Range("N126").Value = -18234.5
$ws.Range("H132").Value = 31312.21
$ws.Range("I132").Value = 6342.148
$ws.Range("K132").Value = 19026.444
$ws.Range("M132").Value = -16496.444
$ws.Range("H136").Value = 73094.836
$ws.Range("J136").Value = 73094.836
$ws.Range("L136").Value = 219284.508
$ws.Range("N136").Value = -224384.508

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7507.4443
$ws.Range("I7").Value = 7429.2666
$ws.Range("K7").Value = 7429.2666
$ws.Range("M7").Value = -7317.2666
$ws.Range("H99").Value = 40950
$ws.Range("I99").Value = 9999
$ws.Range("J99").Value = 56425.5
$ws.Range("K99").Value = 9999
$ws.Range("L99").Value = 56425.5
$ws.Range("M99").Value = -7004
$ws.Range("N99").Value = -62415.5
$ws.Range("H100").Value = 3105.353
$ws.Range("I100").Value = 2286.375
$ws.Range("K100").Value = 2286.375
$ws.Range("M100").Value = -1745.375
$ws.Range("H126").Value = 7507.4443
$ws.Range("I126").Value = 7429.2666
$ws.Range("K126").Value = 22287.7998
$ws.Range("M126").Value = -19817.7998
$ws.Range("H132").Value = 3459.5881
$ws.Range("I132").Value = 2707.25
$ws.Range("K132").Value = 8121.75
$ws.Range("M132").Value = -5591.75
$ws.Range("H136").Value = 150589.97
$ws.Range("I136").Value = 231235.61
$ws.Range("K136").Value = 693706.83
$ws.Range("M136").Value = -691156.83

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H132").Value = 18482.365
$ws.Range("I132").Value = 2225.524
$ws.Range("K132").Value = 6676.572
$ws.Range("M132").Value = -4146.572
$ws.Range("H136").Value = 288631.25
$ws.Range("I136").Value = 315533.53
$ws.Range("K136").Value = 946600.5900000001
$ws.Range("M136").Value = -944050.5900000001

Write-Output "applied changes"